# Split the single "Prerequisites" column into Prerequisites / Corequisites /
# Concurrent / Recommended, inserting three new columns (D, E, F) and pushing
# the old "Terms Typically Offered" column from D out to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert three new blank columns at D..F (old D "Terms Typically Offered"
#    shifts right to G).
$ws.Range("D1:F1").EntireColumn.Insert()

# 2. New header row.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# 3. Default every data row's new Corequisites/Concurrent/Recommended cells to
#    "NA".
$ws.Range("D2:F86").Value = "NA"

# 4. A handful of rows used to bury "Concurrent: ..." / "Recommended: ..."
#    clauses inside the free-text Prerequisites cell. Pull that text out into
#    the new dedicated columns and trim it back out of the Prerequisites text.

# Row 17 - EDUC 427: "Recommended:" clause moves to the Recommended column.
$ws.Range("C17").Value = "Senior standing."
$ws.Range("F17").Value = "ENGL" + [char]160 + "391."
$ws.Range("G17").Value = "F, W, SP "

# Row 29 - EDUC 449: "Concurrent:" clause moves to the Concurrent column.
$ws.Range("C29").Value = "Acceptance into Level I Special Education Credential Program, and completion of all program requirements."
$ws.Range("E29").Value = "EDUC" + [char]160 + "451."
$ws.Range("G29").Value = "SP "

# Row 31 - EDUC 451: "Concurrent:" clause moves to the Concurrent column.
$ws.Range("C31").Value = "Acceptance into Level I Special Education Credential Program; completion of program requirements for the Level I Special Education Program."
$ws.Range("E31").Value = "EDUC" + [char]160 + "449."
$ws.Range("G31").Value = "SP "

# Row 33 - EDUC 455: "Concurrent:" clause moves to the Concurrent column.
$ws.Range("C33").Value = "Acceptance into STEP II or STEP B of the Multiple Subject Teacher Preparation Program."
$ws.Range("E33").Value = "EDUC" + [char]160 + "454."
$ws.Range("G33").Value = "F, W, SP "

# Row 34 - EDUC 456: "Concurrent:" clause moves to the Concurrent column.
$ws.Range("C34").Value = "EDUC" + [char]160 + "454 and EDUC" + [char]160 + "455."
$ws.Range("E34").Value = "EDUC" + [char]160 + "457."
$ws.Range("G34").Value = "F, W, SP "

# Row 35 - EDUC 457: "Concurrent:" clause moves to the Concurrent column.
$ws.Range("C35").Value = "EDUC" + [char]160 + "454 and EDUC" + [char]160 + "455."
$ws.Range("E35").Value = "EDUC" + [char]160 + "456."
$ws.Range("G35").Value = "F, W, SP "
